# Created experiment order generation script
# Re-generates the per-task-order sheets with a freshly rolled timestamp
# suffix for every stimulus filename, rotates which task each physical
# sheet currently holds (GNG <-> TOL), and regenerates the workbook's
# sheet tab names to match the new run ids.

$wb = $excel.ActiveWorkbook

# --- Sheet tab names (rename in place; physical sheet order is unchanged) ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)
$ws4 = $wb.Worksheets.Item(4)
$ws5 = $wb.Worksheets.Item(5)

$ws1.Name = "TOL_TO-16515889474425979"
$ws2.Name = "NB_TO-16515889496411035"
$ws3.Name = "RS_TO-1651588949648028"
$ws4.Name = "GNG_TO-16515889496897817"
$ws5.Name = "vSAT_TO-16515889497662058"

# --- Sheet 1 (was GNG go/GNG_stims order, now becomes the NB MM/ZM order) ---
# grow from 4 data rows (A2:B5) to 6 data rows (A2:B7), copying the
# existing row style down for the two new rows first.
$ws1.Range("A5").Copy($ws1.Range("A6:A7"))

$ws1.Range("B2").Value = "MM_stims-1651588947409778.csv"
$ws1.Range("B3").Value = "ZM_stims-16515889473885555.csv"
$ws1.Range("B4").Value = "MM_stims-16515889474255733.csv"
$ws1.Range("B5").Value = "ZM_stims-16515889474107912.csv"
$ws1.Range("A6").Value = 4
$ws1.Range("B6").Value = "MM_stims-16515889474406004.csv"
$ws1.Range("A7").Value = 5
$ws1.Range("B7").Value = "ZM_stims-16515889474275646.csv"

# --- Sheet 2 (RS-style OB/TB/ZB order; stays 9 data rows) ---
$ws2.Range("B2").Value = "TB-1651588949006439.csv"
$ws2.Range("B3").Value = "ZB-match_2-16515889478013241.csv"
$ws2.Range("B4").Value = "OB-16515889485108843.csv"
$ws2.Range("B5").Value = "TB-16515889492961125.csv"
$ws2.Range("B6").Value = "OB-16515889486048956.csv"
$ws2.Range("B7").Value = "ZB-match_2-1651588948323798.csv"
$ws2.Range("B8").Value = "OB-16515889488953.csv"
$ws2.Range("B9").Value = "ZB-match_4-1651588947709326.csv"
$ws2.Range("B10").Value = "TB-1651588949620282.csv"

# --- Sheet 3 (eyes closed/open resting-state order) is unchanged. ---

# --- Sheet 4 (was NB MM/ZM order, now becomes the GNG go/GNG_stims order) ---
# shrink from 6 data rows (A2:B7) to 4 data rows (A2:B5).
$ws4.Range("B2").Value = "go_stims-16515889496501188.csv"
$ws4.Range("B3").Value = "GNG_stims-16515889496733155.csv"
$ws4.Range("B4").Value = "go_stims-16515889496753514.csv"
$ws4.Range("B5").Value = "GNG_stims-16515889496878226.csv"
$ws4.Range("A6:B7").Delete()

# --- Sheet 5 (SAT/vSAT order; stays 4 data rows) ---
$ws5.Range("B2").Value = "SAT_stims-16515889496946719.csv"
$ws5.Range("B3").Value = "vSAT_stims-16515889497345905.csv"
$ws5.Range("B4").Value = "vSAT_stims-16515889497506745.csv"
$ws5.Range("B5").Value = "SAT_stims-16515889497186453.csv"
